# "Add files via upload" re-save of tweets_supreme_treinamento.xlsx.
# The sheet ("Treinamento") already has columns A (tweet text), B (label)
# and C (a secondary/confidence label that is only populated for some
# rows). This pass updates a batch of previously-entered C values to
# their corrected classification, clears out one stray entry (row 143),
# and leaves the cursor parked on the last cell that was touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Treinamento")

$ws.Range("C3").Value   = 1
$ws.Range("C5").Value   = 3
$ws.Range("C6").Value   = 3
$ws.Range("C11").Value  = 3
$ws.Range("C14").Value  = 1
$ws.Range("C15").Value  = 2
$ws.Range("C16").Value  = 1
$ws.Range("C17").Value  = 4
$ws.Range("C22").Value  = 3
$ws.Range("C26").Value  = 3
$ws.Range("C30").Value  = 3
$ws.Range("C34").Value  = 2
$ws.Range("C40").Value  = 3
$ws.Range("C43").Value  = 2
$ws.Range("C45").Value  = 3
$ws.Range("C46").Value  = 3
$ws.Range("C47").Value  = 3
$ws.Range("C55").Value  = 3
$ws.Range("C69").Value  = 3
$ws.Range("C81").Value  = 3
$ws.Range("C82").Value  = 2
$ws.Range("C84").Value  = 1
$ws.Range("C86").Value  = 3
$ws.Range("C88").Value  = 3
$ws.Range("C90").Value  = 3
$ws.Range("C92").Value  = 1
$ws.Range("C93").Value  = 1
$ws.Range("C94").Value  = 3
$ws.Range("C99").Value  = 1
$ws.Range("C104").Value = 1
$ws.Range("C105").Value = 1
$ws.Range("C107").Value = 2
$ws.Range("C110").Value = 1
$ws.Range("C112").Value = 1
$ws.Range("C113").Value = 1
$ws.Range("C114").Value = 3
$ws.Range("C118").Value = 3
$ws.Range("C123").Value = 3
$ws.Range("C125").Value = 1
$ws.Range("C126").Value = 3
$ws.Range("C127").Value = 2
$ws.Range("C132").Value = 2
$ws.Range("C133").Value = 3
$ws.Range("C135").Value = 3
$ws.Range("C140").Value = 3
$ws.Range("C141").Value = 1
$ws.Range("C144").Value = 1
$ws.Range("C147").Value = 3
$ws.Range("C148").Value = 1
$ws.Range("C155").Value = 3
$ws.Range("C161").Value = 1
$ws.Range("C164").Value = 1
$ws.Range("C166").Value = 3
$ws.Range("C168").Value = 3
$ws.Range("C170").Value = 1
$ws.Range("C174").Value = 1
$ws.Range("C175").Value = 2
$ws.Range("C176").Value = 3
$ws.Range("C185").Value = 3
$ws.Range("C190").Value = 1
$ws.Range("C191").Value = 1
$ws.Range("C201").Value = 1
$ws.Range("C202").Value = 1
$ws.Range("C207").Value = 2
$ws.Range("C210").Value = 1
$ws.Range("C212").Value = 2
$ws.Range("C214").Value = 3
$ws.Range("C218").Value = 2
$ws.Range("C222").Value = 3
$ws.Range("C225").Value = 3
$ws.Range("C226").Value = 2
$ws.Range("C237").Value = 1
$ws.Range("C240").Value = 2
$ws.Range("C241").Value = 2
$ws.Range("C246").Value = 3
$ws.Range("C247").Value = 1
$ws.Range("C249").Value = 3
$ws.Range("C253").Value = 3
$ws.Range("C264").Value = 2
$ws.Range("C271").Value = 1
$ws.Range("C272").Value = 3
$ws.Range("C276").Value = 1
$ws.Range("C285").Value = 3
$ws.Range("C295").Value = 2
$ws.Range("C298").Value = 1

# Row 143 ("PT ..." / label 130): the B value flips from 1 to 0, and the
# stray C143 entry (0) is removed outright rather than just zeroed.
$ws.Range("B143").Value = 0
$ws.Range("C143").ClearContents()

# Column B was left at its default custom width before; re-fit it to its
# contents (matches the "bestFit" auto-sized column in the saved file).
$ws.Columns("B").AutoFit()

# The author's cursor ended on C107 after scrolling the sheet down.
$ws.Application.ActiveWindow.ScrollRow = 106
$ws.Range("C107").Select()
